$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New rows appended at the bottom of the data (rows 193 and 194),
# matching the pattern of the existing columns A:F. Column A holds a
# date-like label that must stay plain text (as the rest of the "Serie"
# column does), so the leading apostrophe forces text entry instead of
# Excel's automatic date conversion; the style is reset back to the
# sheet's default afterwards so no stray formatting is left behind.
$ws.Range("A193").Value = "'06-10-2021"
$ws.Range("B193").Value = 1.69
$ws.Range("C193").Value = 2.54
$ws.Range("D193").Value = 3.27
$ws.Range("E193").Value = 3.96
$ws.Range("F193").Value = -0.12

$ws.Range("A194").Value = "'07-10-2021"
$ws.Range("B194").Value = 1.71
$ws.Range("C194").Value = 2.58
$ws.Range("D194").Value = 3.29
$ws.Range("E194").Value = 4
$ws.Range("F194").Value = -0.11

$ws.Range("A193:A194").Style = "Normal"
